$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.22826733537822008
$ws.Range("B1").Value = 0.22775617450906083
$ws.Range("A2").Value = -0.15706952188144641
$ws.Range("B2").Value = 0.155707033183071
$ws.Range("A3").Value = -0.10599120523820638
$ws.Range("B3").Value = 0.10561410882420752
$ws.Range("A4").Value = -0.09761410891181832
$ws.Range("B4").Value = 0.097282746054345282
$ws.Range("A5").Value = -0.094282746101962189
$ws.Range("B5").Value = 0.093167245180015179
$ws.Range("A6").Value = -0.0096316468721902027
$ws.Range("B6").Value = 0.009574037129183921
$ws.Range("A7").Value = 0.00042596275043838006
$ws.Range("B7").Value = -0.0004296379584949328
$ws.Range("A8").Value = 0.010429637838437866
$ws.Range("B8").Value = -0.01043549905774821
$ws.Range("A9").Value = 0.012435499011243856
$ws.Range("B9").Value = -0.012447659975675318
$ws.Range("A10").Value = 0.014447659931621004
$ws.Range("B10").Value = -0.014447499286580623
$ws.Range("A11").Value = 0.017447499233765207
$ws.Range("B11").Value = -0.017452392335887801
$ws.Range("A12").Value = 0.020952392279024235
$ws.Range("B12").Value = -0.021039435200280732
$ws.Range("A13").Value = 0.024539435148572153
$ws.Range("B13").Value = -0.024616657803123054
$ws.Range("A14").Value = 0.032616657713505859
$ws.Range("B14").Value = -0.032720946734024814
$ws.Range("A15").Value = -0.0080504676448018841
$ws.Range("B15").Value = 0.0080329657944817612
$ws.Range("A16").Value = -0.0060329658282927134
$ws.Range("B16").Value = 0.0060030925594309892
$ws.Range("A17").Value = -0.004003092594011548
$ws.Range("B17").Value = 0.0039999999472879466
$ws.Range("A18").Value = -0.064686652683160872
$ws.Range("B18").Value = 0.064564465661405279
$ws.Range("A19").Value = -0.012091614400012496
$ws.Range("B19").Value = 0.012016390300471969
$ws.Range("A20").Value = -0.0080163903419574467
$ws.Range("B20").Value = 0.0080056622897846097
$ws.Range("A21").Value = -0.0040056623316697682
$ws.Range("B21").Value = 0.0039999999577799983
$ws.Range("A22").Value = -0.045715827992232505
$ws.Range("B22").Value = 0.045501703705172147
$ws.Range("A23").Value = -0.040501703764333818
$ws.Range("B23").Value = 0.040099473700155208
$ws.Range("A24").Value = -0.020099473898324227
$ws.Range("B24").Value = 0.019999999799257928
$ws.Range("A25").Value = -0.081535598357447725
$ws.Range("B25").Value = 0.081448557706512048
$ws.Range("A26").Value = -0.078948557762650751
$ws.Range("B26").Value = 0.078837313479725779
$ws.Range("A27").Value = -0.076337313538812346
$ws.Range("B27").Value = 0.075683117992324966
$ws.Range("A28").Value = -0.073683118059850727
$ws.Range("B28").Value = 0.073243833186950269
$ws.Range("A29").Value = -0.066243833307741973
$ws.Range("B29").Value = 0.06612370101217202
$ws.Range("A30").Value = -0.006123701610680321
$ws.Range("B30").Value = 0.0060671659287847568
$ws.Range("A31").Value = -0.014022487380049853
$ws.Range("B31").Value = 0.014000911027974183
$ws.Range("A32").Value = -0.0040009111813379405
$ws.Range("B32").Value = 0.0039999999001345543

# Column widths: target OOXML width attrs are 15.7109375 (col A) and 16.42578125 (col B).
# This runtime quantizes ColumnWidth to a 6px-per-char grid (+5px padding) when serializing,
# so we pick the ColumnWidth that lands closest to each target after that quantization.
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
